$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.86201166666667
$ws.Range("H2").Value = 32.586035
$ws.Range("I2").Value = 0.08282714153498995
$ws.Range("J2").Value = 0.08282714153498993
$ws.Range("M2").Value = 13.788265
$ws.Range("N2").Value = 41.364795
$ws.Range("O2").Value = 0.7604887294515349
$ws.Range("P2").Value = 0.760488729451535
$ws.Range("Q2").Value = 149.7682952930917
$ws.Range("R2").Value = 1347.914657637825
$ws.Range("S2").Value = 0.06298910763004696
$ws.Range("T2").Value = 0.06298910763004696
$ws.Range("G3").Value = 10.86201166666667
$ws.Range("H3").Value = 32.586035
$ws.Range("I3").Value = 0.08282714153498995
$ws.Range("J3").Value = 0.08282714153498993
$ws.Range("O3").Value = 0.1168722076973129
$ws.Range("P3").Value = 0.1168722076973129
$ws.Range("Q3").Value = 23.01645065350333
$ws.Range("R3").Value = 207.14805588153
$ws.Range("S3").Value = 0.009680190888452078
$ws.Range("T3").Value = 0.009680190888452076
$ws.Range("G4").Value = 10.86201166666667
$ws.Range("H4").Value = 32.586035
$ws.Range("I4").Value = 0.08282714153498995
$ws.Range("J4").Value = 0.08282714153498993
$ws.Range("M4").Value = 1.925545
$ws.Range("N4").Value = 5.776635
$ws.Range("O4").Value = 0.1062030118040055
$ws.Range("P4").Value = 0.1062030118040055
$ws.Range("Q4").Value = 20.91529225469166
$ws.Range("R4").Value = 188.237630292225
$ws.Range("S4").Value = 0.00879649189013257
$ws.Range("T4").Value = 0.00879649189013257
$ws.Range("G5").Value = 10.86201166666667
$ws.Range("H5").Value = 32.586035
$ws.Range("I5").Value = 0.08282714153498995
$ws.Range("J5").Value = 0.08282714153498993
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2979986666666667
$ws.Range("N5").Value = 0.893996
$ws.Range("O5").Value = 0.0164360510471466
$ws.Range("P5").Value = 0.0164360510471466
$ws.Range("Q5").Value = 3.236864993984444
$ws.Range("R5").Value = 29.13178494586
$ws.Range("S5").Value = 0.001361351126358331
$ws.Range("T5").Value = 0.001361351126358331
$ws.Range("G6").Value = 92.04504633333333
$ws.Range("I6").Value = 0.7018799384686454
$ws.Range("J6").Value = 0.7018799384686453
$ws.Range("M6").Value = 13.788265
$ws.Range("N6").Value = 41.364795
$ws.Range("O6").Value = 0.7604887294515349
$ws.Range("P6").Value = 0.760488729451535
$ws.Range("Q6").Value = 1269.141490781278
$ws.Range("R6").Value = 11422.27341703151
$ws.Range("S6").Value = 0.5337717826335416
$ws.Range("T6").Value = 0.5337717826335416
$ws.Range("G7").Value = 92.04504633333333
$ws.Range("I7").Value = 0.7018799384686454
$ws.Range("J7").Value = 0.7018799384686453
$ws.Range("O7").Value = 0.1168722076973129
$ws.Range("P7").Value = 0.1168722076973129
$ws.Range("S7").Value = 0.08203025794728472
$ws.Range("T7").Value = 0.08203025794728472
$ws.Range("G8").Value = 92.04504633333333
$ws.Range("I8").Value = 0.7018799384686454
$ws.Range("J8").Value = 0.7018799384686453
$ws.Range("M8").Value = 1.925545
$ws.Range("N8").Value = 5.776635
$ws.Range("O8").Value = 0.1062030118040055
$ws.Range("P8").Value = 0.1062030118040055
$ws.Range("Q8").Value = 177.2368787419183
$ws.Range("R8").Value = 1595.131908677265
$ws.Range("S8").Value = 0.07454176339018018
$ws.Range("T8").Value = 0.07454176339018019
$ws.Range("G9").Value = 92.04504633333333
$ws.Range("I9").Value = 0.7018799384686454
$ws.Range("J9").Value = 0.7018799384686453
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2979986666666667
$ws.Range("N9").Value = 0.893996
$ws.Range("O9").Value = 0.0164360510471466
$ws.Range("P9").Value = 0.0164360510471466
$ws.Range("Q9").Value = 27.42930108060489
$ws.Range("R9").Value = 246.863709725444
$ws.Range("S9").Value = 0.01153613449763877
$ws.Range("T9").Value = 0.01153613449763877
$ws.Range("G10").Value = 10.015157
$ws.Range("H10").Value = 30.045471
$ws.Range("I10").Value = 0.0763695392520887
$ws.Range("J10").Value = 0.07636953925208868
$ws.Range("M10").Value = 13.788265
$ws.Range("N10").Value = 41.364795
$ws.Range("O10").Value = 0.7604887294515349
$ws.Range("P10").Value = 0.760488729451535
$ws.Range("Q10").Value = 138.091638732605
$ws.Range("R10").Value = 1242.824748593445
$ws.Range("S10").Value = 0.05807817387462005
$ws.Range("T10").Value = 0.05807817387462005
$ws.Range("G11").Value = 10.015157
$ws.Range("H11").Value = 30.045471
$ws.Range("I11").Value = 0.0763695392520887
$ws.Range("J11").Value = 0.07636953925208868
$ws.Range("O11").Value = 0.1168722076973129
$ws.Range("P11").Value = 0.1168722076973129
$ws.Range("Q11").Value = 21.221977470802
$ws.Range("R11").Value = 190.997797237218
$ws.Range("S11").Value = 0.0089254766532182
$ws.Range("T11").Value = 0.0089254766532182
$ws.Range("G12").Value = 10.015157
$ws.Range("H12").Value = 30.045471
$ws.Range("I12").Value = 0.0763695392520887
$ws.Range("J12").Value = 0.07636953925208868
$ws.Range("M12").Value = 1.925545
$ws.Range("N12").Value = 5.776635
$ws.Range("O12").Value = 0.1062030118040055
$ws.Range("P12").Value = 0.1062030118040055
$ws.Range("Q12").Value = 19.284635485565
$ws.Range("R12").Value = 173.561719370085
$ws.Range("S12").Value = 0.008110675078656036
$ws.Range("T12").Value = 0.008110675078656036
$ws.Range("G13").Value = 10.015157
$ws.Range("H13").Value = 30.045471
$ws.Range("I13").Value = 0.0763695392520887
$ws.Range("J13").Value = 0.07636953925208868
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.2979986666666667
$ws.Range("N13").Value = 0.893996
$ws.Range("O13").Value = 0.0164360510471466
$ws.Range("P13").Value = 0.0164360510471466
$ws.Range("Q13").Value = 2.984503432457334
$ws.Range("R13").Value = 26.860530892116
$ws.Range("S13").Value = 0.001255213645594396
$ws.Range("T13").Value = 0.001255213645594395
$ws.Range("G14").Value = 18.21851333333333
$ws.Range("H14").Value = 54.65554
$ws.Range("I14").Value = 0.138923380744276
$ws.Range("J14").Value = 0.138923380744276
$ws.Range("M14").Value = 13.788265
$ws.Range("N14").Value = 41.364795
$ws.Range("O14").Value = 0.7604887294515349
$ws.Range("P14").Value = 0.760488729451535
$ws.Range("Q14").Value = 251.2016897460334
$ws.Range("R14").Value = 2260.8152077143
$ws.Range("S14").Value = 0.1056496653133263
$ws.Range("T14").Value = 0.1056496653133263
$ws.Range("G15").Value = 18.21851333333333
$ws.Range("H15").Value = 54.65554
$ws.Range("I15").Value = 0.138923380744276
$ws.Range("J15").Value = 0.138923380744276
$ws.Range("O15").Value = 0.1168722076973129
$ws.Range("P15").Value = 0.1168722076973129
$ws.Range("Q15").Value = 38.60477469414667
$ws.Range("R15").Value = 347.44297224732
$ws.Range("S15").Value = 0.01623628220835791
$ws.Range("T15").Value = 0.01623628220835791
$ws.Range("G16").Value = 18.21851333333333
$ws.Range("H16").Value = 54.65554
$ws.Range("I16").Value = 0.138923380744276
$ws.Range("J16").Value = 0.138923380744276
$ws.Range("M16").Value = 1.925545
$ws.Range("N16").Value = 5.776635
$ws.Range("O16").Value = 0.1062030118040055
$ws.Range("P16").Value = 0.1062030118040055
$ws.Range("Q16").Value = 35.08056725643333
$ws.Range("R16").Value = 315.7251053079
$ws.Range("S16").Value = 0.0147540814450367
$ws.Range("T16").Value = 0.0147540814450367
$ws.Range("G17").Value = 18.21851333333333
$ws.Range("H17").Value = 54.65554
$ws.Range("I17").Value = 0.138923380744276
$ws.Range("J17").Value = 0.138923380744276
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2979986666666667
$ws.Range("N17").Value = 0.893996
$ws.Range("O17").Value = 0.0164360510471466
$ws.Range("P17").Value = 0.0164360510471466
$ws.Range("Q17").Value = 5.429092681982223
$ws.Range("R17").Value = 48.86183413784
$ws.Range("S17").Value = 0.002283351777555104
$ws.Range("T17").Value = 0.002283351777555104
